$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rtn4"
$ws.Range("C2").Value = "Lingo1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 40.797777
$ws.Range("H2").Value = 122.393331
$ws.Range("I2").Value = 0.2689231481273683
$ws.Range("J2").Value = 0.2689231481273683
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9237416666666666
$ws.Range("N2").Value = 2.771225
$ws.Range("O2").Value = 0.9776944302049534
$ws.Range("P2").Value = 0.9776944302049534
$ws.Range("Q2").Value = 37.686606522275
$ws.Range("R2").Value = 339.179458700475
$ws.Range("S2").Value = 0.2629246640773096
$ws.Range("T2").Value = 0.2629246640773097

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rtn4"
$ws.Range("C3").Value = "Lingo1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 40.797777
$ws.Range("H3").Value = 122.393331
$ws.Range("I3").Value = 0.2689231481273683
$ws.Range("J3").Value = 0.2689231481273683
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.02107466666666667
$ws.Range("N3").Value = 0.063224
$ws.Range("O3").Value = 0.02230556979504659
$ws.Range("P3").Value = 0.02230556979504659
$ws.Range("Q3").Value = 0.8597995510160001
$ws.Range("R3").Value = 7.738195959144002
$ws.Range("S3").Value = 0.005998484050058666
$ws.Range("T3").Value = 0.005998484050058667

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rtn4"
$ws.Range("C4").Value = "Lingo1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 46.219831
$ws.Range("H4").Value = 138.659493
$ws.Range("I4").Value = 0.3046632285488233
$ws.Range("J4").Value = 0.3046632285488233
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9237416666666666
$ws.Range("N4").Value = 2.771225
$ws.Range("O4").Value = 0.9776944302049534
$ws.Range("P4").Value = 0.9776944302049534
$ws.Range("Q4").Value = 42.69518372099166
$ws.Range("R4").Value = 384.256653488925
$ws.Range("S4").Value = 0.2978675416404433
$ws.Range("T4").Value = 0.2978675416404433

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rtn4"
$ws.Range("C5").Value = "Lingo1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 46.219831
$ws.Range("H5").Value = 138.659493
$ws.Range("I5").Value = 0.3046632285488233
$ws.Range("J5").Value = 0.3046632285488233
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.02107466666666667
$ws.Range("N5").Value = 0.063224
$ws.Range("O5").Value = 0.02230556979504659
$ws.Range("P5").Value = 0.02230556979504659
$ws.Range("Q5").Value = 0.9740675317146668
$ws.Range("R5").Value = 8.766607785432001
$ws.Range("S5").Value = 0.006795686908380009
$ws.Range("T5").Value = 0.006795686908380009

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Rtn4"
$ws.Range("C6").Value = "Lingo1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 64.69033266666666
$ws.Range("H6").Value = 194.070998
$ws.Range("I6").Value = 0.4264136233238083
$ws.Range("J6").Value = 0.4264136233238083
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.9237416666666666
$ws.Range("N6").Value = 2.771225
$ws.Range("O6").Value = 0.9776944302049534
$ws.Range("P6").Value = 0.9776944302049534
$ws.Range("Q6").Value = 59.75715571472777
$ws.Range("R6").Value = 537.81440143255
$ws.Range("S6").Value = 0.4169022244872004
$ws.Range("T6").Value = 0.4169022244872004

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Rtn4"
$ws.Range("C7").Value = "Lingo1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 64.69033266666666
$ws.Range("H7").Value = 194.070998
$ws.Range("I7").Value = 0.4264136233238083
$ws.Range("J7").Value = 0.4264136233238083
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02107466666666667
$ws.Range("N7").Value = 0.063224
$ws.Range("O7").Value = 0.02230556979504659
$ws.Range("P7").Value = 0.02230556979504659
$ws.Range("Q7").Value = 1.363327197505778
$ws.Range("R7").Value = 12.269944777552
$ws.Range("S7").Value = 0.009511398836607912
$ws.Range("T7").Value = 0.009511398836607912

